# Update "想去人数" (want-to-go count) figures across the four sheets of the
# workbook, matching the refreshed data snapshot from the site regeneration.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 501
$ws.Range("F6").Value = 933
$ws.Range("F7").Value = 470
$ws.Range("F9").Value = 2221
$ws.Range("F10").Value = 635
$ws.Range("F11").Value = 298
$ws.Range("F14").Value = 188
$ws.Range("F15").Value = 2236
$ws.Range("F16").Value = 689
$ws.Range("F17").Value = 13879
$ws.Range("F18").Value = 7
$ws.Range("F19").Value = 1301
$ws.Range("F20").Value = 63
$ws.Range("F21").Value = 567
$ws.Range("F22").Value = 139
$ws.Range("F23").Value = 36
$ws.Range("F24").Value = 149
$ws.Range("F25").Value = 97
$ws.Range("F27").Value = 280
$ws.Range("F29").Value = 3
$ws.Range("F31").Value = 31

# --- Sheet: 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 21
$ws.Range("F9").Value = 151
$ws.Range("F11").Value = 86
$ws.Range("F12").Value = 64
$ws.Range("F18").Value = 28

# --- Sheet: 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5726
$ws.Range("F4").Value = 475

# --- Sheet: 全部类型 (All Types, combined view) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 475
$ws.Range("F7").Value = 501
$ws.Range("F8").Value = 933
$ws.Range("F10").Value = 470
$ws.Range("F11").Value = 21
$ws.Range("F12").Value = 2221
$ws.Range("F13").Value = 635
$ws.Range("F14").Value = 298
$ws.Range("F20").Value = 188
$ws.Range("F21").Value = 151
$ws.Range("F23").Value = 2236
$ws.Range("F24").Value = 689
$ws.Range("F25").Value = 86
$ws.Range("F26").Value = 64
$ws.Range("F27").Value = 1301
$ws.Range("F28").Value = 63
$ws.Range("F29").Value = 567
$ws.Range("F30").Value = 139
$ws.Range("F31").Value = 36
$ws.Range("F32").Value = 149
$ws.Range("F33").Value = 97
$ws.Range("F38").Value = 280
$ws.Range("F49").Value = 31
